$wb = $excel.ActiveWorkbook

# --- Sheet 1 (Vmax): change selection to header row A1:H1 ---
$wsVmax = $wb.Worksheets.Item("Vmax")
$wsVmax.Range("A1:H1").Select() | Out-Null

# --- Sheet 2 (Km): change selection to single cell F10 (also clears tabSelected) ---
$wsKm = $wb.Worksheets.Item("Km")
$wsKm.Range("F10").Select() | Out-Null

# --- Add new sheet "litterChemistry" after Km, at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLitter = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLitter.Name = "litterChemistry"

# --- Populate the litter chemistry ANOVA results table ---
$data = @(
    @("functionalGroup","timePoint","Vegetation","Precipitation","timePoint x Precipitation","timePoint x Vegetation","Vegetation x Precipitation","Three-way"),
    @("glycosidicBond","o","***","**","o","o","*","o"),
    @("C_O_stretching","*","***","**","o","o","o","o"),
    @("carboEster","o","***","*","o","o","o","o"),
    @("lipid","o","***","o","o","o","*","o"),
    @("alkane","**","***","o","o","***","*","o"),
    @("amide","o","o","***","o","o","o","o")
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt $data[$r].Count; $c++) {
        $wsLitter.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Bold header row
$wsLitter.Range("A1:H1").Font.Bold = $true

# Column widths to match source table (values chosen so the engine's
# internal pixel-quantized ColumnWidth storage rounds to the closest
# achievable width to the original authored file's bestFit widths)
$wsLitter.Columns.Item(1).ColumnWidth = 14.0
$wsLitter.Columns.Item(2).ColumnWidth = 8.333333333333332
$wsLitter.Columns.Item(3).ColumnWidth = 9.333333333333332
$wsLitter.Columns.Item(4).ColumnWidth = 10.833333333333332
$wsLitter.Columns.Item(5).ColumnWidth = 21.166666666666668
$wsLitter.Columns.Item(6).ColumnWidth = 19.666666666666668
$wsLitter.Columns.Item(7).ColumnWidth = 22.166666666666668
$wsLitter.Columns.Item(8).ColumnWidth = 9.0

# Select the full data range, matching the saved selection state
$wsLitter.Range("A1:H7").Select() | Out-Null
